$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header "no of cuts" in F1, matching style of other header cells (B1:E1)
$ws.Range("F1").Value = "no of cuts"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update row 2 text values
$ws.Range("C2").Value = "منطقة 2"
$ws.Range("D2").Value = "ثقب"
$ws.Range("E2").Value = "مثقاب"

# New numeric value in F2
$ws.Range("F2").Value = 3
